# Weekly data refresh: insert the latest "Escarola" (Primera / Segunda)
# price observations dated 2021-10-07 at the top of the data block (row 401),
# pushing all existing rows down by two and extending the used range from
# R484 to R486.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 401 down by inserting two fresh rows.
$ws.Rows("401:402").Insert()

$newRows = @(
    @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", "2021-10-07", 15, 100112033, "Lechuga", "Escarola", "Primera", 120, 4500, 5000, 4750, "`$/caja 12 unidades", "Región de Arica y Parinacota", 396, 12, "Hortaliza"),
    @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", "2021-10-07", 15, 100112033, "Lechuga", "Escarola", "Segunda", 120, 4500, 5000, 4750, "`$/caja 18 unidades", "Región de Arica y Parinacota", 264, 18, "Hortaliza")
)

$startRow = 401
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
